# Auto-generated Excel COM-interop edit script reproducing the target diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: Normality_Results -----------------------------------------
# Recomputed normality-test statistics (tiny floating point jitter from a
# re-run of the analysis) across columns B,C,D,E,F,G,H,I.
$wsNorm = $wb.Worksheets.Item("Normality_Results")
$wsNorm.Range("F2").Value = 0.89169317483901978
$wsNorm.Range("G2").Value = 0.028893610462546349
$wsNorm.Range("B3").Value = 0.92305392026901245
$wsNorm.Range("C3").Value = 0.24323886632919309
$wsNorm.Range("D4").Value = 0.88559812307357788
$wsNorm.Range("E4").Value = 0.069872625172138214
$wsNorm.Range("F5").Value = 0.85639059543609619
$wsNorm.Range("G5").Value = 0.006833657156676054
$wsNorm.Range("B6").Value = 0.84110546112060547
$wsNorm.Range("C6").Value = 0.016855573281645771
$wsNorm.Range("H7").Value = 0.84225279092788696
$wsNorm.Range("I7").Value = 0.0039639859460294247
$wsNorm.Range("F9").Value = 0.82710939645767212
$wsNorm.Range("G9").Value = 0.002256682375445962
$wsNorm.Range("D10").Value = 0.81390684843063354
$wsNorm.Range("E10").Value = 0.0074545401148498058
$wsNorm.Range("F10").Value = 0.86510264873504639
$wsNorm.Range("G10").Value = 0.0096479486674070358
$wsNorm.Range("B12").Value = 0.86701583862304688
$wsNorm.Range("C12").Value = 0.038109641522169113
$wsNorm.Range("F12").Value = 0.85951763391494751
$wsNorm.Range("G12").Value = 0.0077278297394514084
$wsNorm.Range("B13").Value = 0.92419636249542236
$wsNorm.Range("C13").Value = 0.25255483388900762
$wsNorm.Range("D13").Value = 0.84203445911407471
$wsNorm.Range("E13").Value = 0.017344687134027481
$wsNorm.Range("F13").Value = 0.78874766826629639
$wsNorm.Range("G13").Value = 0.00058990460820496082
$wsNorm.Range("F14").Value = 0.80984175205230713
$wsNorm.Range("G14").Value = 0.001215857220813632
$wsNorm.Range("F15").Value = 0.89029496908187866
$wsNorm.Range("G15").Value = 0.02723001129925251
$wsNorm.Range("I15").Value = 0.0095877842977643013
$wsNorm.Range("F16").Value = 0.74556714296340942
$wsNorm.Range("G16").Value = 0.00014813755115028471
$wsNorm.Range("D18").Value = 0.85528254508972168
$wsNorm.Range("E18").Value = 0.026219135150313381
$wsNorm.Range("F19").Value = 0.87606239318847656
$wsNorm.Range("G19").Value = 0.015041562728583809
$wsNorm.Range("H19").Value = 0.87546300888061523
$wsNorm.Range("I19").Value = 0.014676352962851519
$wsNorm.Range("D20").Value = 0.7444191575050354
$wsNorm.Range("E20").Value = 0.0011095058871433141

# --- Sheet: Wilcoxon_A_Results -----------------------------------------
# Recomputed Wilcoxon p-values (column C) + apply a 3-decimal number format.
$wsWA = $wb.Worksheets.Item("Wilcoxon_A_Results")
$wsWA.Range("C2").Value = 0.1520926324815797
$wsWA.Range("C3").Value = 0.78252792474006738
$wsWA.Range("C4").Value = 0.40762594770278088
$wsWA.Range("C5").Value = 0.0033680111449411539
$wsWA.Range("C6").Value = 0.011866216879384949
$wsWA.Range("C7").Value = 0.89023005494349738
$wsWA.Range("C8").Value = 0.60490715489275892
$wsWA.Range("C9").Value = 0.66461162954100628
$wsWA.Range("C10").Value = 0.58621368107313998
$wsWA.Range("C11").Value = 0.89023005494349738
$wsWA.Range("C12").Value = 0.62442260383854364
$wsWA.Range("C13").Value = 0.031424346521782533
$wsWA.Range("C14").Value = 0.088973011701813334
$wsWA.Range("C15").Value = 0.83021759177786336
$wsWA.Range("C16").Value = 0.076730336987448799
$wsWA.Range("C17").Value = 0.669028074498903
$wsWA.Range("C18").Value = 0.76559448399576402
$wsWA.Range("C19").Value = 1
$wsWA.Range("C20").Value = 0.77282999268444752
$wsWA.Range("C2:C20").NumberFormat = "0.000"
$wsWA.Range("F15").Select()

# --- Sheet: Wilcoxon_B_Results -----------------------------------------
# Recomputed Wilcoxon p-values (column C) + apply a 3-decimal number format.
$wsWB = $wb.Worksheets.Item("Wilcoxon_B_Results")
$wsWB.Range("C2").Value = 0.067172958725004298
$wsWB.Range("C3").Value = 0.19655330119230929
$wsWB.Range("C4").Value = 0.53272590644282347
$wsWB.Range("C5").Value = 0.002281937253315448
$wsWB.Range("C6").Value = 0.7589513020672245
$wsWB.Range("C7").Value = 0.8025873486341526
$wsWB.Range("C8").Value = 0.066362066456708876
$wsWB.Range("C9").Value = 0.26611935562031019
$wsWB.Range("C10").Value = 0.19252671718116021
$wsWB.Range("C11").Value = 0.63642497301956835
$wsWB.Range("C12").Value = 0.68355287527265496
$wsWB.Range("C13").Value = 0.023536373622090701
$wsWB.Range("C14").Value = 0.43740644290235942
$wsWB.Range("C15").Value = 0.2364984548439383
$wsWB.Range("C16").Value = 0.0079101778825628969
$wsWB.Range("C17").Value = 0.74458688845732768
$wsWB.Range("C18").Value = 0.029873144019229029
$wsWB.Range("C19").Value = 0.24295250731663051
$wsWB.Range("C20").Value = 0.52419397534258827
$wsWB.Range("C2:C20").NumberFormat = "0.000"
$wsWB.Range("E16").Select()

# --- Sheet: Mann_Whitney_Non_AI_Results ---------------------------------
# Values unchanged; apply a 3-decimal number format to the p-value column.
$wsMWN = $wb.Worksheets.Item("Mann_Whitney_Non_AI_Results")
$wsMWN.Range("C2:C20").NumberFormat = "0.000"

# --- Sheet: Mann_Whitney_AI_Results --------------------------------------
# Values unchanged; apply a 3-decimal number format to the p-value column.
$wsMWAI = $wb.Worksheets.Item("Mann_Whitney_AI_Results")
$wsMWAI.Range("C2:C20").NumberFormat = "0.000"
$wsMWAI.Range("F16").Select()

# --- Final view state ----------------------------------------------------
# Mann_Whitney_Non_AI_Results ends up the active/selected sheet & cell.
$wsMWN.Activate()
$wsMWN.Range("F21").Select()

